# "more work on presentations"
#
# 1) The footer "datetimeFigureOut" date field cached on the slide master
#    and on every slide layout gets bumped from 11/5/20 to 11/12/20.
# 2) The picture named "Picture 4" on slide 3 is nudged up a bit
#    (its vertical offset changes from 2171700 EMU to 2135124 EMU).

$p = $ppt.ActivePresentation

$oldDate = "11/5/20"
$newDate = "11/12/20"

# --- Slide master: "Date Placeholder" shape ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout under the master: "Date Placeholder" shape ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 3: move "Picture 4" up slightly ---
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $sh = $slide3.Shapes.Item($i)
    if ($sh.Name -eq "Picture 4") {
        # 2135124 EMU -> points, nudged so the stored EMU value round-trips exactly
        $sh.Top = 168.12001
    }
}
